$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '65.804.33'
$ws.Range('E2').Value = '  -4.73%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.526.27'
$ws.Range('E3').Value = '  -5.22%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.994'
$ws.Range('E4').Value = '  -0.44%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '577.08'
$ws.Range('E5').Value = '  -6.33%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '183.35'
$ws.Range('E6').Value = '  +0.33%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.517.93'
$ws.Range('E7').Value = '  -5.31%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.601'
$ws.Range('E8').Value = '  -4.99%  '
$ws.Range('E9').Value = '  +0.03%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.660'
$ws.Range('E10').Value = '  -8.61%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '53.55'
$ws.Range('E11').Value = '  -7.00%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.141'
$ws.Range('E12').Value = '  -12.75%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000247'
$ws.Range('E13').Value = '  -16.46%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '9.68'
$ws.Range('E14').Value = '  -9.82%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.077.02'
$ws.Range('E15').Value = '  -5.37%  '
$ws.Range('E16').Value = '  -0.84%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.490.92'
$ws.Range('E17').Value = '  -6.20%  '
$ws.Range('B18').Value = 'Chainlink'
$ws.Range('C18').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '18.10'
$ws.Range('E18').Value = '  -7.21%  '
$ws.Range('B19').Value = 'WrappedBTC'
$ws.Range('C19').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '65.215.47'
$ws.Range('E19').Value = '  -5.19%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.00'
$ws.Range('E20').Value = '  -7.72%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.04'
$ws.Range('E21').Value = '  -8.18%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '388.62'
$ws.Range('E22').Value = '  -6.43%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.23'
$ws.Range('E23').Value = '  -10.66%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '83.71'
$ws.Range('E24').Value = '  -6.21%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.83'
$ws.Range('E25').Value = '  -7.35%  '
$ws.Range('B26').Value = 'InternetComputer(DFINITY)'
$ws.Range('C26').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '12.18'
$ws.Range('E26').Value = '  -4.68%  '
$ws.Range('B27').Value = 'LEO'
$ws.Range('C27').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '5.97'
$ws.Range('E27').Value = '  -1.84%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.20'
$ws.Range('E28').Value = '  -7.17%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '3.54'
$ws.Range('E29').Value = '  -9.31%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.80'
$ws.Range('E30').Value = '  -8.66%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '30.59'
$ws.Range('E31').Value = '  -7.60%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.69'
$ws.Range('E32').Value = '  -9.12%  '
$ws.Range('B33').Value = 'Cosmos'
$ws.Range('C33').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '11.87'
$ws.Range('E33').Value = '  -5.46%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '64.97'
$ws.Range('E34').Value = '  -1.00%  '
$ws.Range('B35').Value = 'Bittensor'
$ws.Range('C35').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '606.13'
$ws.Range('E35').Value = '  -0.72%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.111'
$ws.Range('E36').Value = '  -7.82%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '41.07'
$ws.Range('E37').Value = '  -6.80%  '
$ws.Range('E38').Value = '  +0.15%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.993'
$ws.Range('E39').Value = '  -0.59%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.368'
$ws.Range('E40').Value = '  -9.76%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0₃0733'
$ws.Range('E41').Value = '  -17.83%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.129'
$ws.Range('E42').Value = '  -7.38%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.839.90'
$ws.Range('E43').Value = '  +1.56%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.74'
$ws.Range('E44').Value = '  -10.45%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0401'
$ws.Range('E45').Value = '  -9.46%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.39'
$ws.Range('E46').Value = '  -10.35%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.129'
$ws.Range('E47').Value = '  -5.25%  '
$ws.Range('B48').Value = 'Monero'
$ws.Range('C48').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '137.13'
$ws.Range('E48').Value = '  -2.58%  '
$ws.Range('B49').Value = 'ApeXProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.94'
$ws.Range('E49').Value = '  -4.96%  '
$ws.Range('B50').Value = 'THORChain'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '8.14'
$ws.Range('E50').Value = '  -12.16%  '
$ws.Range('B51').Value = 'WEMIXToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.43'
$ws.Range('E51').Value = '  -10.75%  '
